$d = $word.ActiveDocument

# Step 1: The "Assigned To: ... {{assignedto}}" and "Created By: ... {{createdby}}"
# lines currently live in one paragraph, separated by a manual line break (vertical
# tab char 11 in Range.Text / <w:br/> in XML). Replace that line break with a real
# paragraph mark so the line break becomes a paragraph split.
$d.Content.Find.Execute([string][char]11, $false, $false, $false, $false, $false, `
                         $true, 1, $false, [string][char]13, 2)

# Step 2: Split the bold "Assigned To: " run into "Assigned To" (now red) and ": ".
# Find "Assigned To" within the (now separate) paragraph and color it red;
# Word will automatically split the run so the trailing ": " keeps the old formatting.
$r = $d.Content
$r.Find.Execute("Assigned To", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Font.Color = 255

# Step 3: Move the "_GoBack" bookmark from right after "{{description}}" to the very
# start of the new "Created By" paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$p3 = $d.Paragraphs.Item(3)
$start = $p3.Range.Start
$newBookmarkRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
